$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.436.12'
$ws.Range('E2').Value = '  -1.66%  '
$ws.Range('D3').Value = '2.633.35'
$ws.Range('E3').Value = '  -0.72%  '
$ws.Range('D5').Value = "'581.83"
$ws.Range('E5').Value = '  -2.50%  '
$ws.Range('D6').Value = "'156.76"
$ws.Range('E6').Value = '  +0.80%  '
$ws.Range('E7').Value = '  +3.11%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  -3.41%  '
$ws.Range('D10').Value = "'5.81"
$ws.Range('E10').Value = '  +0.34%  '
$ws.Range('E11').Value = '  -1.31%  '
$ws.Range('E12').Value = '  -0.12%  '
$ws.Range('D13').Value = "'28.59"
$ws.Range('E13').Value = '  -0.61%  '
$ws.Range('E14').Value = '  -4.76%  '
$ws.Range('D15').Value = '3.109.47'
$ws.Range('E15').Value = '  -0.65%  '
$ws.Range('D16').Value = '64.240.26'
$ws.Range('E16').Value = '  -1.79%  '
$ws.Range('D17').Value = '2.626.94'
$ws.Range('E17').Value = '  -2.75%  '
$ws.Range('D18').Value = "'12.25"
$ws.Range('E18').Value = '  -2.96%  '
$ws.Range('E19').Value = '  -1.76%  '
$ws.Range('D20').Value = "'7.44"
$ws.Range('E20').Value = '  -0.34%  '
$ws.Range('D21').Value = "'346.71"
$ws.Range('E21').Value = '  -0.55%  '
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('D23').Value = "'68.28"
$ws.Range('E23').Value = '  -0.90%  '
$ws.Range('D24').Value = "'1.77"
$ws.Range('E24').Value = '  +7.03%  '
$ws.Range('E25').Value = '  +0.22%  '
$ws.Range('D26').Value = "'9.44"
$ws.Range('E26').Value = '  -1.91%  '
$ws.Range('D27').Value = "'588.96"
$ws.Range('E27').Value = '  +9.41%  '
$ws.Range('D28').Value = "'1.59"
$ws.Range('E28').Value = '  +0.35%  '
$ws.Range('D29').Value = "'7.98"
$ws.Range('E29').Value = '  +0.40%  '
$ws.Range('E30').Value = '  -1.61%  '
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('E32').Value = '  -1.38%  '
$ws.Range('D33').Value = "'6.67"
$ws.Range('E33').Value = '  +4.20%  '
$ws.Range('D34').Value = "'1.72"
$ws.Range('E34').Value = '  -1.68%  '
$ws.Range('D35').Value = "'5.34"
$ws.Range('E35').Value = '  -1.63%  '
$ws.Range('D36').Value = "'0.412"
$ws.Range('E36').Value = '  -1.67%  '
$ws.Range('D37').Value = "'20.03"
$ws.Range('E37').Value = '  -1.66%  '
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('E39').Value = '  +1.00%  '
$ws.Range('D40').Value = "'154.94"
$ws.Range('E40').Value = '  -0.52%  '
$ws.Range('E42').Value = '  +6.04%  '
$ws.Range('D43').Value = "'157.96"
$ws.Range('E43').Value = '  -1.93%  '
$ws.Range('E44').Value = '  -1.55%  '
$ws.Range('D45').Value = "'23.29"
$ws.Range('E45').Value = '  +3.57%  '
$ws.Range('D46').Value = "'0.0602"
$ws.Range('E46').Value = '  -0.59%  '
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('E48').Value = '  +2.70%  '
$ws.Range('E49').Value = '  -0.42%  '
$ws.Range('D50').Value = "'19.21"
$ws.Range('E50').Value = '  -2.08%  '
$ws.Range('E51').Value = '  -6.05%  '
